$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 623.8588
$ws.Cells.Item(17, 10).Value = 623.8588
$ws.Cells.Item(17, 12).Value = 1871.5764
$ws.Cells.Item(17, 14).Value = -2207.5764
$ws.Cells.Item(43, 8).Value = 1099
$ws.Cells.Item(43, 10).Value = 1099
$ws.Cells.Item(43, 12).Value = 1099
$ws.Cells.Item(43, 14).Value = -1237
$ws.Cells.Item(51, 8).Value = 3993.75
$ws.Cells.Item(51, 10).Value = 4000
$ws.Cells.Item(51, 12).Value = 4000
$ws.Cells.Item(51, 14).Value = -4968
$ws.Cells.Item(92, 8).Value = 807
$ws.Cells.Item(92, 9).Value = 811.7143
$ws.Cells.Item(92, 11).Value = 811.7143
$ws.Cells.Item(92, 13).Value = 436.2857
$ws.Cells.Item(107, 8).Value = 178.77777
$ws.Cells.Item(107, 9).Value = 203.57143
$ws.Cells.Item(107, 10).Value = 92
$ws.Cells.Item(107, 11).Value = 203.57143
$ws.Cells.Item(107, 12).Value = 92
$ws.Cells.Item(107, 13).Value = 1716.42857
$ws.Cells.Item(107, 14).Value = -3932
$ws.Cells.Item(110, 8).Value = 50193
$ws.Cells.Item(110, 10).Value = 50193
$ws.Cells.Item(110, 12).Value = 50193
$ws.Cells.Item(110, 14).Value = -58373
$ws.Cells.Item(116, 8).Value = 8247.833000000001
$ws.Cells.Item(116, 9).Value = 8560
$ws.Cells.Item(116, 10).Value = 8024.857
$ws.Cells.Item(116, 11).Value = 8560
$ws.Cells.Item(116, 12).Value = 8024.857
$ws.Cells.Item(116, 13).Value = -5118
$ws.Cells.Item(116, 14).Value = -14908.857
$ws.Cells.Item(117, 8).Value = 99999
$ws.Cells.Item(117, 10).Value = 99999
$ws.Cells.Item(117, 12).Value = 99999
$ws.Cells.Item(117, 14).Value = -109177
$ws.Cells.Item(132, 8).Value = 1419.9767
$ws.Cells.Item(132, 9).Value = 1501.7838
$ws.Cells.Item(132, 10).Value = 915.5
$ws.Cells.Item(132, 11).Value = 4505.3514
$ws.Cells.Item(132, 12).Value = 2746.5
$ws.Cells.Item(132, 13).Value = -1975.3514
$ws.Cells.Item(132, 14).Value = -7806.5
$ws.Cells.Item(133, 8).Value = 77179.914
$ws.Cells.Item(133, 10).Value = 77179.914
$ws.Cells.Item(133, 12).Value = 77179.914
$ws.Cells.Item(133, 14).Value = -87299.914
$ws.Cells.Item(134, 8).Value = 99999
$ws.Cells.Item(134, 10).Value = 99999
$ws.Cells.Item(134, 12).Value = 99999
$ws.Cells.Item(134, 14).Value = -110139
$ws.Cells.Item(135, 8).Value = 2011.3636
$ws.Cells.Item(135, 9).Value = 1833.5
$ws.Cells.Item(135, 10).Value = 2811.75
$ws.Cells.Item(135, 11).Value = 16501.5
$ws.Cells.Item(135, 12).Value = 25305.75
$ws.Cells.Item(135, 13).Value = -13966.5
$ws.Cells.Item(135, 14).Value = -30375.75
$ws.Cells.Item(136, 8).Value = 96495.836
$ws.Cells.Item(136, 10).Value = 96495.836
$ws.Cells.Item(136, 12).Value = 96495.836
$ws.Cells.Item(136, 14).Value = -106695.836
$ws.Cells.Item(139, 8).Value = 99999
$ws.Cells.Item(139, 10).Value = 99999
$ws.Cells.Item(139, 12).Value = 99999
$ws.Cells.Item(139, 14).Value = -110279
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2261.125
$ws.Cells.Item(61, 10).Value = 2227.2727
$ws.Cells.Item(61, 12).Value = 2227.2727
$ws.Cells.Item(61, 14).Value = -2651.2727
$ws.Cells.Item(108, 8).Value = 88996.8
$ws.Cells.Item(108, 10).Value = 88996.8
$ws.Cells.Item(108, 12).Value = 88996.8
$ws.Cells.Item(108, 14).Value = -96676.8
$ws.Cells.Item(136, 8).Value = 2261.125
$ws.Cells.Item(136, 10).Value = 2227.2727
$ws.Cells.Item(136, 12).Value = 6681.8181
$ws.Cells.Item(136, 14).Value = -11781.8181
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(108, 8).Value = 99995
$ws.Cells.Item(108, 10).Value = 99995
$ws.Cells.Item(108, 12).Value = 99995
$ws.Cells.Item(108, 14).Value = -107675
$ws.Cells.Item(109, 8).Value = 73139.14
$ws.Cells.Item(109, 10).Value = 73139.14
$ws.Cells.Item(109, 12).Value = 73139.14
$ws.Cells.Item(109, 14).Value = -75913.14
$ws.Cells.Item(132, 8).Value = 96282.28999999999
$ws.Cells.Item(132, 10).Value = 96282.28999999999
$ws.Cells.Item(132, 12).Value = 96282.28999999999
$ws.Cells.Item(132, 14).Value = -106402.29
$ws.Cells.Item(138, 8).Value = 99999
$ws.Cells.Item(138, 10).Value = 99999
$ws.Cells.Item(138, 12).Value = 99999
$ws.Cells.Item(138, 14).Value = -110279
$ws.Cells.Item(140, 8).Value = 94017.89
$ws.Cells.Item(140, 10).Value = 43412.332
$ws.Cells.Item(140, 12).Value = 43412.332
$ws.Cells.Item(140, 14).Value = -53772.332
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(105, 8).Value = 3914.1667
$ws.Cells.Item(105, 9).Value = 871.25
$ws.Cells.Item(105, 11).Value = 871.25
$ws.Cells.Item(105, 13).Value = 875.75
$ws.Cells.Item(114, 8).Value = 36116.625
$ws.Cells.Item(114, 10).Value = 36116.625
$ws.Cells.Item(114, 12).Value = 36116.625
$ws.Cells.Item(114, 14).Value = -44794.625
$ws.Cells.Item(118, 8).Value = 99999
$ws.Cells.Item(118, 10).Value = 99999
$ws.Cells.Item(118, 12).Value = 99999
$ws.Cells.Item(118, 14).Value = -103313
$ws.Cells.Item(120, 8).Value = 49987.5
$ws.Cells.Item(120, 10).Value = 49987.5
$ws.Cells.Item(120, 12).Value = 49987.5
$ws.Cells.Item(120, 14).Value = -57245.5
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 1366.8572
$ws.Cells.Item(46, 9).Value = 120
$ws.Cells.Item(46, 10).Value = 1574.6666
$ws.Cells.Item(46, 11).Value = 360
$ws.Cells.Item(46, 12).Value = 4723.9998
$ws.Cells.Item(46, 13).Value = -269
$ws.Cells.Item(46, 14).Value = -4905.9998
$ws.Cells.Item(86, 8).Value = 703.3333
$ws.Cells.Item(86, 10).Value = 760.625
$ws.Cells.Item(86, 12).Value = 2281.875
$ws.Cells.Item(86, 14).Value = -4653.875
$ws.Cells.Item(89, 8).Value = 703.3333
$ws.Cells.Item(89, 10).Value = 760.625
$ws.Cells.Item(89, 12).Value = 6845.625
$ws.Cells.Item(89, 14).Value = -18701.625
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(116, 8).Value = 59163.332
$ws.Cells.Item(116, 10).Value = 59163.332
$ws.Cells.Item(116, 12).Value = 59163.332
$ws.Cells.Item(116, 14).Value = -68341.33199999999
$ws.Cells.Item(122, 8).Value = 296706.22
$ws.Cells.Item(122, 9).Value = 386657.78
$ws.Cells.Item(122, 11).Value = 1159973.34
$ws.Cells.Item(122, 13).Value = -1157523.34
$ws.Cells.Item(135, 8).Value = 52276.617
$ws.Cells.Item(135, 10).Value = 52276.617
$ws.Cells.Item(135, 12).Value = 52276.617
$ws.Cells.Item(135, 14).Value = -62416.617
$ws.Cells.Item(140, 8).Value = 97331.11
$ws.Cells.Item(140, 10).Value = 98185
$ws.Cells.Item(140, 12).Value = 98185
$ws.Cells.Item(140, 14).Value = -108545
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 7388.5
$ws.Cells.Item(46, 9).Value = 18091.334
$ws.Cells.Item(46, 11).Value = 18091.334
$ws.Cells.Item(46, 13).Value = -17903.334
$ws.Cells.Item(118, 8).Value = 50290.91
$ws.Cells.Item(118, 10).Value = 51320
$ws.Cells.Item(118, 12).Value = 51320
$ws.Cells.Item(118, 14).Value = -54634
$ws.Cells.Item(122, 8).Value = 75004430
$ws.Cells.Item(122, 9).Value = 111115304
$ws.Cells.Item(122, 10).Value = 28576168
$ws.Cells.Item(122, 11).Value = 333345912
$ws.Cells.Item(122, 12).Value = 85728504
$ws.Cells.Item(122, 13).Value = -333343462
$ws.Cells.Item(122, 14).Value = -85733404
$ws.Cells.Item(123, 8).Value = 75421.75
$ws.Cells.Item(123, 10).Value = 78997.71000000001
$ws.Cells.Item(123, 12).Value = 78997.71000000001
$ws.Cells.Item(123, 14).Value = -88797.71000000001
$ws.Cells.Item(129, 8).Value = 95673.625
$ws.Cells.Item(129, 9).Value = 63195
$ws.Cells.Item(129, 11).Value = 63195
$ws.Cells.Item(129, 13).Value = -58195
$ws.Cells.Item(136, 8).Value = 1912.16
$ws.Cells.Item(136, 9).Value = 1510
$ws.Cells.Item(136, 10).Value = 2138.375
$ws.Cells.Item(136, 11).Value = 4530
$ws.Cells.Item(136, 12).Value = 6415.125
$ws.Cells.Item(136, 13).Value = -1980
$ws.Cells.Item(136, 14).Value = -11515.125
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1796.8857
$ws.Cells.Item(132, 9).Value = 1524.8846
$ws.Cells.Item(132, 10).Value = 2582.6667
$ws.Cells.Item(132, 11).Value = 4574.6538
$ws.Cells.Item(132, 12).Value = 7748.000100000001
$ws.Cells.Item(132, 13).Value = -2044.6538
$ws.Cells.Item(132, 14).Value = -12808.0001
